$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write row data (rows 2-12) ---
$ws.Range("A2").Value = '2026-02-02 12:54:56'
$ws.Range("B2").Value = '【完全在宅】AI×Web開発エンジニア募集!業務自動化・AI機能開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5483480'
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = '🔥AI,Ai ◆開発,自動化'

$ws.Range("A3").Value = '2026-02-02 12:54:56'
$ws.Range("B3").Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5483966'
$ws.Range("G3").Value = 225
$ws.Range("H3").Value = '🔥Next.js ◆開発 ◇アプリ'

$ws.Range("A4").Value = '2026-02-02 12:54:56'
$ws.Range("B4").Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5483967'
$ws.Range("G4").Value = 218
$ws.Range("H4").Value = '🔥Next.js ◆開発 ◇アプリ'

$ws.Range("A5").Value = '2026-02-02 12:54:56'
$ws.Range("B5").Value = '【急募】楽天RPP広告自動化ツールの開発依頼'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5484001'
$ws.Range("G5").Value = 213
$ws.Range("H5").Value = '◆ツール,開発'

$ws.Range("A6").Value = '2026-02-02 12:54:56'
$ws.Range("B6").Value = 'X(旧twitter)のロック解除自動化システム構築'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5484231'
$ws.Range("G6").Value = 103
$ws.Range("H6").Value = '◆自動化'

$ws.Range("A7").Value = '2026-02-02 12:54:56'
$ws.Range("B7").Value = '【急募】新しいWebサービスの開発パートナーを探しています!'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5483482'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'

$ws.Range("A8").Value = '2026-02-02 12:54:56'
$ws.Range("B8").Value = '四柱推命の命式自動計算プログラム(Web/Excel)の開発依頼'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5484177'
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = '◆開発'

$ws.Range("A9").Value = '2026-02-02 12:54:56'
$ws.Range("B9").Value = '美容皮膚科向け LINE連携型BtoB SaaS(MVP) の開発案件'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5483503'
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = '◆開発'

$ws.Range("A10").Value = '2026-02-02 12:54:56'
$ws.Range("B10").Value = '【急募】Notion×Slackでのオンライン講座運営システム構築'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5483854'
$ws.Range("G10").Value = 28

$ws.Range("A11").Value = '2026-02-02 12:54:56'
$ws.Range("B11").Value = '【1,000件以上の案件からマッチ】中長期で参画可能なエンジニアを大募集|気軽に相談OK'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5484020'
$ws.Range("G11").Value = 25

$ws.Range("A12").Value = '2026-02-02 12:54:56'
$ws.Range("B12").Value = '【市場調査】海外向けデジタルサービスの価値評価依頼'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5483504'
$ws.Range("G12").Value = 13

# --- Rebuild hyperlinks for column F (rows 2-12) in order so rIds come out sequential ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5483480')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5483966')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5483967')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5484001')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5484231')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5483482')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5484177')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5483503')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5483854')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5484020')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5483504')

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 45.1
$ws.Columns.Item(8).ColumnWidth = 18.1

